# Apply odds updates for Jogos_da_Semana_FlashScore_2025-01-30.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = 1.48
$ws.Range("Q4").Value = 1.77
$ws.Range("R4").Value = 1.97
$ws.Range("AB4").Value = 1.63
$ws.Range("I5").Value = 3.1
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 1.94
$ws.Range("R5").Value = 1.79
$ws.Range("U5").Value = 4.2
$ws.Range("V5").Value = 1.23
$ws.Range("AB5").Value = 1.63
$ws.Range("AE5").Value = 10
$ws.Range("AA6").Value = 1.58
$ws.Range("S7").Value = 1.85
$ws.Range("T7").Value = 2
$ws.Range("G8").Value = 2.5
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 3.1
$ws.Range("L8").Value = 3.25
$ws.Range("AC8").Value = 9.5
$ws.Range("AD8").Value = 13
$ws.Range("AE8").Value = 10
$ws.Range("AP8").Value = 10
$ws.Range("AQ8").Value = 26
$ws.Range("G9").Value = 3.05
$ws.Range("H9").Value = 3.1
$ws.Range("I9").Value = 2.3
$ws.Range("J9").Value = 3.55
$ws.Range("K9").Value = 2.02
$ws.Range("L9").Value = 2.9
$ws.Range("O9").Value = 1.28
$ws.Range("P9").Value = 3.05
$ws.Range("S9").Value = 1.82
$ws.Range("T9").Value = 1.78
$ws.Range("W9").Value = 2.87
$ws.Range("X9").Value = 1.31
$ws.Range("Y9").Value = 1.4
$ws.Range("Z9").Value = 2.52
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 17
$ws.Range("AE9").Value = 10.5
$ws.Range("AF9").Value = 40
$ws.Range("AG9").Value = 26
$ws.Range("AH9").Value = 30
$ws.Range("AI9").Value = 9.5
$ws.Range("AJ9").Value = 6
$ws.Range("AL9").Value = 50
$ws.Range("AO9").Value = 11.75
$ws.Range("AP9").Value = 8.75
$ws.Range("AQ9").Value = 24
$ws.Range("AR9").Value = 18.5
$ws.Range("AS9").Value = 26
$ws.Range("G10").Value = 1.77
$ws.Range("G11").Value = 1.92
$ws.Range("I11").Value = 3.9
$ws.Range("J11").Value = 2.88
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 1.53
$ws.Range("P11").Value = 2.38
$ws.Range("Q11").Value = 2.03
$ws.Range("R11").Value = 1.78
$ws.Range("S11").Value = 2.7
$ws.Range("T11").Value = 1.44
$ws.Range("Y11").Value = 1.62
$ws.Range("Z11").Value = 2.2
$ws.Range("AA11").Value = 2.25
$ws.Range("AB11").Value = 1.57
$ws.Range("AD11").Value = 8
$ws.Range("AF11").Value = 17
$ws.Range("AI11").Value = 6
$ws.Range("AO11").Value = 19
$ws.Range("G14").Value = 3.1
$ws.Range("H14").Value = 3.65
$ws.Range("I14").Value = 2.07
$ws.Range("J14").Value = 3.4
$ws.Range("K14").Value = 2.35
$ws.Range("L14").Value = 2.55
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 9.5
$ws.Range("O14").Value = 1.15
$ws.Range("P14").Value = 4.75
$ws.Range("S14").Value = 1.47
$ws.Range("T14").Value = 2.52
$ws.Range("W14").Value = 2.1
$ws.Range("X14").Value = 1.65
$ws.Range("Y14").Value = 1.26
$ws.Range("Z14").Value = 3.45
$ws.Range("AA14").Value = 1.42
$ws.Range("AB14").Value = 2.65
$ws.Range("AC14").Value = 15.5
$ws.Range("AD14").Value = 22
$ws.Range("AE14").Value = 11.25
$ws.Range("AF14").Value = 45
$ws.Range("AG14").Value = 22
$ws.Range("AH14").Value = 22
$ws.Range("AI14").Value = 9.5
$ws.Range("AJ14").Value = 7.9
$ws.Range("AK14").Value = 10.75
$ws.Range("AL14").Value = 32
$ws.Range("AM14").Value = 150
$ws.Range("AN14").Value = 12
$ws.Range("AO14").Value = 13.5
$ws.Range("AP14").Value = 8.75
$ws.Range("AQ14").Value = 22
$ws.Range("AR14").Value = 14
$ws.Range("AS14").Value = 17.5
$ws.Range("I16").Value = 1.27
$ws.Range("L16").Value = 1.69
$ws.Range("G17").Value = 2
$ws.Range("I17").Value = 3.5
$ws.Range("J17").Value = 2.62
$ws.Range("L17").Value = 3.75
$ws.Range("S17").Value = 1.88
$ws.Range("T17").Value = 1.93
$ws.Range("W17").Value = 3.25
$ws.Range("X17").Value = 1.33
$ws.Range("AG17").Value = 17
$ws.Range("I18").Value = 1.19
$ws.Range("L18").Value = 1.54
$ws.Range("G19").Value = 13.5
$ws.Range("H19").Value = 6.4
$ws.Range("I19").Value = 1.18
$ws.Range("J19").Value = 9.5
$ws.Range("P19").Value = 6.1
$ws.Range("S19").Value = 1.3
$ws.Range("T19").Value = 3.2
$ws.Range("X19").Value = 2
$ws.Range("AA19").Value = 1.75
$ws.Range("AB19").Value = 1.95
$ws.Range("AC19").Value = 50
$ws.Range("AD19").Value = 120
$ws.Range("AE19").Value = 40
$ws.Range("AF19").Value = 500
$ws.Range("AI19").Value = 26
$ws.Range("AJ19").Value = 14.5
$ws.Range("AK19").Value = 23
$ws.Range("AO19").Value = 8.25
$ws.Range("AP19").Value = 10
$ws.Range("AQ19").Value = 8.25
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 4.45
$ws.Range("J20").Value = 4.7
$ws.Range("K20").Value = 2.57
$ws.Range("O20").Value = 1.13
$ws.Range("P20").Value = 5.2
$ws.Range("T20").Value = 2.7
$ws.Range("W20").Value = 2
$ws.Range("X20").Value = 1.72
$ws.Range("Y20").Value = 1.24
$ws.Range("Z20").Value = 3.7
$ws.Range("AC20").Value = 23
$ws.Range("AD20").Value = 37
$ws.Range("AE20").Value = 16.5
$ws.Range("AF20").Value = 90
$ws.Range("AG20").Value = 40
$ws.Range("AJ20").Value = 9.5
$ws.Range("AN20").Value = 11.25
$ws.Range("AO20").Value = 9.75
$ws.Range("AQ20").Value = 12.5
$ws.Range("G21").Value = 1.6
$ws.Range("H21").Value = 4.1
$ws.Range("I21").Value = 4.7
$ws.Range("J21").Value = 2.07
$ws.Range("L21").Value = 4.6
$ws.Range("Y21").Value = 1.27
$ws.Range("Z21").Value = 3.4
$ws.Range("AA21").Value = 1.6
$ws.Range("AB21").Value = 2.22
$ws.Range("AH21").Value = 19.5
$ws.Range("AJ21").Value = 8.25
$ws.Range("AM21").Value = 250
$ws.Range("AN21").Value = 18
$ws.Range("AO21").Value = 32
$ws.Range("AP21").Value = 15
$ws.Range("AQ21").Value = 80
